# Swap the data of duplicate-location observation rows that were
# mismatched in the source export. Columns that already hold identical
# values between the two rows of a pair are swapped too (harmless,
# since swapping equal values is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <-> Row 5 : Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor
$pair1Row1 = 3
$pair1Row2 = 5
$pair1Cols = @("A","B","E","F","G","H")
foreach ($col in $pair1Cols) {
    $cell1 = $ws.Range("$col" + "$pair1Row1")
    $cell2 = $ws.Range("$col" + "$pair1Row2")
    $tmp = $cell1.Value2
    $cell1.Value2 = $cell2.Value2
    $cell2.Value2 = $tmp
}

# Row 7 <-> Row 8 : Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor
$pair2Row1 = 7
$pair2Row2 = 8
$pair2Cols = @("A","B","E","F","G","H")
foreach ($col in $pair2Cols) {
    $cell1 = $ws.Range("$col" + "$pair2Row1")
    $cell2 = $ws.Range("$col" + "$pair2Row2")
    $tmp = $cell1.Value2
    $cell1.Value2 = $cell2.Value2
    $cell2.Value2 = $tmp
}

# Row 9 <-> Row 10 : Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
#                    Vetenskapligt namn, Auktor, Ost, Nord, Starttid, Sluttid
$pair3Row1 = 9
$pair3Row2 = 10
$pair3Cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")
foreach ($col in $pair3Cols) {
    $cell1 = $ws.Range("$col" + "$pair3Row1")
    $cell2 = $ws.Range("$col" + "$pair3Row2")
    $tmp = $cell1.Value2
    $cell1.Value2 = $cell2.Value2
    $cell2.Value2 = $tmp
}
